$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# --- ERM table view: add the missing "foreign" key markers in column F
#     for rows 4-6 (mirrors the existing F3 "foreign" label: italic font,
#     right-aligned) ---
$ws.Range("F4").Value = "foreign"
$ws.Range("F4").Font.Italic = $true
$ws.Range("F4").HorizontalAlignment = -4152   # xlRight

$ws.Range("F5").Value = "foreign"
$ws.Range("F5").Font.Italic = $true
$ws.Range("F5").HorizontalAlignment = -4152   # xlRight

$ws.Range("F6").Value = "foreign"
$ws.Range("F6").Font.Italic = $true
$ws.Range("F6").HorizontalAlignment = -4152   # xlRight

# --- Move the viewport/selection as the author left it: scrolled so
#     column D is the left-most visible column, with L12 selected ---
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L12").Select()
